$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 207
$ws1.Range("F7").Value = 53
$ws1.Range("F10").Value = 4476
$ws1.Range("F11").Value = 4476
$ws1.Range("F13").Value = 439
$ws1.Range("F16").Value = 3912
$ws1.Range("F19").Value = 31
$ws1.Range("F20").Value = 153
$ws1.Range("F21").Value = 3292
$ws1.Range("F25").Value = 2763
$ws1.Range("F27").Value = 112
$ws1.Range("F29").Value = 122
$ws1.Range("F30").Value = 161
$ws1.Range("F32").Value = 65
$ws1.Range("F33").Value = 37
$ws1.Range("F34").Value = 21
$ws1.Range("F35").Value = 48
$ws1.Range("F37").Value = 5051
$ws1.Range("F42").Value = 8
$ws1.Range("F43").Value = 982
$ws1.Range("F44").Value = 373
$ws1.Range("F46").Value = 1863
$ws1.Range("F49").Value = 666

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 8
$ws2.Range("F8").Value = 52

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 8
$ws4.Range("F8").Value = 207
$ws4.Range("F11").Value = 53
$ws4.Range("F14").Value = 4476
$ws4.Range("F15").Value = 4476
$ws4.Range("F17").Value = 52
$ws4.Range("F19").Value = 439
$ws4.Range("F22").Value = 3912
$ws4.Range("F25").Value = 3292
$ws4.Range("F26").Value = 2763
$ws4.Range("F28").Value = 112
$ws4.Range("F29").Value = 122
$ws4.Range("F30").Value = 161
$ws4.Range("F32").Value = 65
$ws4.Range("F33").Value = 21
$ws4.Range("F35").Value = 48
$ws4.Range("F39").Value = 5051
$ws4.Range("F45").Value = 982
$ws4.Range("F46").Value = 373
$ws4.Range("F48").Value = 1863
$ws4.Range("F50").Value = 666
